$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 176, shifting the existing rows 176:213 down to 177:214
$ws.Rows.Item(176).EntireRow.Insert()

# Populate the newly inserted row 176 with the new weekly record
$ws.Cells.Item(176, 1).Value = 5
$ws.Cells.Item(176, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(176, 3).Value = "Maule"
$ws.Cells.Item(176, 4).Value = 44641
$ws.Cells.Item(176, 5).Value = 7
$ws.Cells.Item(176, 6).Value = 100112009
$ws.Cells.Item(176, 7).Value = "Acelga"
$ws.Cells.Item(176, 8).Value = "Sin especificar"
$ws.Cells.Item(176, 9).Value = "Primera"
$ws.Cells.Item(176, 10).Value = 400
$ws.Cells.Item(176, 11).Value = 3500
$ws.Cells.Item(176, 12).Value = 3500
$ws.Cells.Item(176, 13).Value = 3500
$ws.Cells.Item(176, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(176, 15).Value = "Región del Maule"
$ws.Cells.Item(176, 16).Value = 875
$ws.Cells.Item(176, 17).Value = 4
$ws.Cells.Item(176, 18).Value = "Hortaliza"
